$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:G51 to text format so numeric-looking values (percentages,
# plain numbers) are stored as literal text, matching the source data.
$ws.Range("D2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "303.23"
$ws.Range("E2").Value = "5.37%"
$ws.Range("G2").Value = "11"

$ws.Range("D3").Value = "34.95"
$ws.Range("E3").Value = "12.90%"
$ws.Range("G3").Value = "11"

$ws.Range("D4").Value = "5.174"
$ws.Range("E4").Value = "5.24%"
$ws.Range("G4").Value = "11"

$ws.Range("E5").Value = "6.03%"
$ws.Range("G5").Value = "11"

$ws.Range("D6").Value = "2.319"
$ws.Range("E6").Value = "-0.45%"
$ws.Range("G6").Value = "11"

$ws.Range("D7").Value = "8.026"
$ws.Range("E7").Value = "3.67%"
$ws.Range("G7").Value = "11"

$ws.Range("D8").Value = "4.002"
$ws.Range("E8").Value = "7.54%"
$ws.Range("G8").Value = "11"

$ws.Range("D9").Value = "0.9273"
$ws.Range("E9").Value = "2.66%"
$ws.Range("G9").Value = "11"

$ws.Range("E10").Value = "11.43%"
$ws.Range("G10").Value = "11"

$ws.Range("D11").Value = "0.1825"
$ws.Range("E11").Value = "8.13%"
$ws.Range("G11").Value = "11"

$ws.Range("E12").Value = "5.54%"
$ws.Range("G12").Value = "11"

$ws.Range("D13").Value = "0.03478"
$ws.Range("E13").Value = "11.35%"
$ws.Range("G13").Value = "11"

$ws.Range("D14").Value = "0.09895"
$ws.Range("E14").Value = "-0.39%"
$ws.Range("G14").Value = "11"

$ws.Range("D15").Value = "0.001483"
$ws.Range("E15").Value = "-0.79%"
$ws.Range("G15").Value = "11"

$ws.Range("D16").Value = "0.04620"
$ws.Range("E16").Value = "2.54%"
$ws.Range("G16").Value = "11"

$ws.Range("D17").Value = "0.005806"
$ws.Range("E17").Value = "0.58%"
$ws.Range("G17").Value = "11"

$ws.Range("D18").Value = "3.480"
$ws.Range("E18").Value = "-0.47%"
$ws.Range("G18").Value = "11"

$ws.Range("D19").Value = "2.105"
$ws.Range("E19").Value = "0.35%"
$ws.Range("G19").Value = "11"

$ws.Range("D20").Value = "0.3439"
$ws.Range("E20").Value = "3.28%"
$ws.Range("G20").Value = "11"

$ws.Range("E21").Value = "2.63%"
$ws.Range("G21").Value = "11"

$ws.Range("D22").Value = "4.603"
$ws.Range("E22").Value = "8.86%"
$ws.Range("G22").Value = "11"

$ws.Range("D23").Value = "0.2341"
$ws.Range("E23").Value = "11.47%"
$ws.Range("G23").Value = "11"

$ws.Range("D24").Value = "0.001223"
$ws.Range("E24").Value = "0.98%"
$ws.Range("G24").Value = "11"

$ws.Range("D25").Value = "0.004425"
$ws.Range("E25").Value = "6.22%"
$ws.Range("G25").Value = "11"

$ws.Range("D26").Value = "0.0001304"
$ws.Range("E26").Value = "0.31%"
$ws.Range("G26").Value = "11"

$ws.Range("D27").Value = "0.0003420"
$ws.Range("E27").Value = "0.78%"
$ws.Range("G27").Value = "11"

$ws.Range("G28").Value = "11"

$ws.Range("G29").Value = "11"

$ws.Range("G30").Value = "11"

$ws.Range("G31").Value = "11"

$ws.Range("G32").Value = "11"

$ws.Range("G33").Value = "11"

$ws.Range("G34").Value = "11"

$ws.Range("G35").Value = "11"

$ws.Range("G36").Value = "11"

$ws.Range("G37").Value = "11"

$ws.Range("G38").Value = "11"

$ws.Range("D39").Value = "0.01765"
$ws.Range("E39").Value = "12.10%"
$ws.Range("G39").Value = "11"

$ws.Range("D40").Value = "0.04705"
$ws.Range("E40").Value = "6.00%"
$ws.Range("G40").Value = "11"

$ws.Range("D41").Value = "0.007637"
$ws.Range("E41").Value = "3.78%"
$ws.Range("G41").Value = "11"

$ws.Range("D42").Value = "0.1406"
$ws.Range("E42").Value = "5.93%"
$ws.Range("G42").Value = "11"

$ws.Range("D43").Value = "0.007078"
$ws.Range("E43").Value = "-25.28%"
$ws.Range("G43").Value = "11"

$ws.Range("D44").Value = "0.002209"
$ws.Range("E44").Value = "-0.93%"
$ws.Range("G44").Value = "11"

$ws.Range("D45").Value = "0.009212"
$ws.Range("E45").Value = "2.53%"
$ws.Range("G45").Value = "11"

$ws.Range("D46").Value = "0.00005925"
$ws.Range("E46").Value = "-3.10%"
$ws.Range("G46").Value = "11"

$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.29%"
$ws.Range("G47").Value = "11"

$ws.Range("E48").Value = "11.15%"
$ws.Range("G48").Value = "11"

$ws.Range("D49").Value = "0.002704"
$ws.Range("E49").Value = "35.12%"
$ws.Range("G49").Value = "11"

$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").Value = "0.29%"
$ws.Range("G50").Value = "11"

$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.29%"
$ws.Range("G51").Value = "11"

# Clear the temporary text-number-format override so cells fall back to
# the workbook default style (matches original formatting, no explicit s=).
$ws.Range("D2:G51").Style = "Normal"
